$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52

# Copy formatting of an existing similar row (row 33) into the new row first,
# so fonts / borders / row height etc. come along for the ride.
$ws.Range("A33:I33").Copy()
$ws.Range("A52:I52").PasteSpecial(-4122)
$ws.Rows.Item(52).RowHeight = $ws.Rows.Item(14).RowHeight

# --- Column I : Keywords2 (created first so it lands at the expected shared-string index) ---
$params = "coyni_mobile.tests.CustomerProfileTest,`ntestLogInSessions,`n-ppinHeading,`n-ppin,`n-pexpHeading,`n-pdescription,`n-ppassword,`n-psubHeading,`n-pnewPassword,`n-pConfirmPassword,`n-psucessHeading,`n-psucessDesc,`n-pvalidateChangePassword,`n-pendSessnSucessHeading"
$ws.Range("I52").Value = $params

# --- Column G : tags ---
$ws.Range("G52").Value = "Profile-Login Sessions"

# --- Column A : Test Case Name ---
$ws.Range("A52").Value = "Verify Login Sessions in Profile"
$ws.Range("A52").HorizontalAlignment = 1

# --- Column B : Execute ---
$ws.Range("B52").Value = "No"

# --- Column C : Test Data Path ---
$ws.Range("C52").Value = "testdata_3_0_customer.xls,profile"

# --- Column D : Iteration Mode ---
$ws.Range("D29").Copy()
$ws.Range("D52").PasteSpecial(-4122)
$ws.Range("D52").Value = "RunOneIteration"
$ws.Range("D52").VerticalAlignment = -4108
$ws.Range("D52").HorizontalAlignment = 1

# --- Column E : Start Iteration ---
$ws.Range("D29").Copy()
$ws.Range("E52").PasteSpecial(-4122)
$ws.Range("E52").VerticalAlignment = -4108
$ws.Range("E52").HorizontalAlignment = 1
$ws.Range("E52").Value = "'1"

# --- Column F : Stop Iteration ---
$ws.Range("D29").Copy()
$ws.Range("F52").PasteSpecial(-4122)
$ws.Range("F52").VerticalAlignment = -4108
$ws.Range("F52").HorizontalAlignment = 1
$ws.Range("F52").Value = "'3"

# --- Column H : Keywords1 ---
$ws.Range("H12").Copy()
$ws.Range("H52").PasteSpecial(-4122)
$ws.Range("H52").Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-password,`n-ppin"

$ws.Range("A1").Select()
Write-Host "done"
